$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 7639.2
$ws.Range("I2").Value = 1781.5714
$ws.Range("K2").Value = 1781.5714
$ws.Range("M2").Value = -1668.5714
$ws.Range("H28").Value = 5270.8667
$ws.Range("I28").Value = 484.1
$ws.Range("K28").Value = 484.1
$ws.Range("M28").Value = 0.8999999999999773
$ws.Range("H55").Value = 122.63636
$ws.Range("J55").Value = 166
$ws.Range("L55").Value = 166
$ws.Range("N55").Value = -594
$ws.Range("H62").Value = 76924500
$ws.Range("I62").Value = 90910470
$ws.Range("J62").Value = 1602
$ws.Range("K62").Value = 90910470
$ws.Range("L62").Value = 1602
$ws.Range("M62").Value = -90909846
$ws.Range("N62").Value = -2850
$ws.Range("H65").Value = 76924500
$ws.Range("I65").Value = 90910470
$ws.Range("J65").Value = 1602
$ws.Range("K65").Value = 454552350
$ws.Range("L65").Value = 8010
$ws.Range("M65").Value = -454549230
$ws.Range("N65").Value = -14250
$ws.Range("H98").Value = 13823
$ws.Range("I98").Value = 15845.363
$ws.Range("K98").Value = 15845.363
$ws.Range("M98").Value = -14347.363
$ws.Range("H112").Value = 2620.4827
$ws.Range("I112").Value = 199
$ws.Range("J112").Value = 2706.9644
$ws.Range("K112").Value = 597
$ws.Range("L112").Value = 8120.8932
$ws.Range("M112").Value = 511
$ws.Range("N112").Value = -10336.8932
$ws.Range("H122").Value = 13823
$ws.Range("I122").Value = 15845.363
$ws.Range("K122").Value = 47536.089
$ws.Range("M122").Value = -45086.089
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4280883
$ws.Range("I32").Value = 5293101
$ws.Range("J32").Value = 7074.222
$ws.Range("K32").Value = 5293101
$ws.Range("L32").Value = 7074.222
$ws.Range("M32").Value = -5292814
$ws.Range("N32").Value = -7648.222
$ws.Range("H45").Value = 4436
$ws.Range("I45").Value = 4436
$ws.Range("K45").Value = 4436
$ws.Range("M45").Value = -4059
$ws.Range("H54").Value = 62000
$ws.Range("J54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("N54").Value = -41538
$ws.Range("H61").Value = 7223.7856
$ws.Range("I61").Value = 5068.8887
$ws.Range("J61").Value = 11102.6
$ws.Range("K61").Value = 5068.8887
$ws.Range("L61").Value = 11102.6
$ws.Range("M61").Value = -4856.8887
$ws.Range("N61").Value = -11526.6
$ws.Range("H97").Value = 1962131
$ws.Range("I97").Value = 2647841.2
$ws.Range("K97").Value = 2647841.2
$ws.Range("M97").Value = -2647345.2
$ws.Range("H132").Value = 10877.272
$ws.Range("I132").Value = 5400
$ws.Range("J132").Value = 15441.667
$ws.Range("K132").Value = 16200
$ws.Range("L132").Value = 46325.001
$ws.Range("M132").Value = -13670
$ws.Range("N132").Value = -51385.001
$ws.Range("H136").Value = 7223.7856
$ws.Range("I136").Value = 5068.8887
$ws.Range("J136").Value = 11102.6
$ws.Range("K136").Value = 15206.6661
$ws.Range("L136").Value = 33307.8
$ws.Range("M136").Value = -12656.6661
$ws.Range("N136").Value = -38407.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2697.139
$ws.Range("J107").Value = 11137.5
$ws.Range("L107").Value = 11137.5
$ws.Range("N107").Value = -14977.5
$ws.Range("H132").Value = 100911.8
$ws.Range("J132").Value = 100911.8
$ws.Range("L132").Value = 100911.8
$ws.Range("N132").Value = -111031.8
$ws.Range("H134").Value = 6621.8486
$ws.Range("I134").Value = 2856.6667
$ws.Range("K134").Value = 8570.000100000001
$ws.Range("M134").Value = -6035.000100000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4788.607
$ws.Range("I58").Value = 2977.2104
$ws.Range("K58").Value = 2977.2104
$ws.Range("M58").Value = -2774.2104
$ws.Range("H132").Value = 62089.285
$ws.Range("I132").Value = 3792.2
$ws.Range("K132").Value = 11376.6
$ws.Range("M132").Value = -8846.599999999999
$ws.Range("H136").Value = 4788.607
$ws.Range("I136").Value = 2977.2104
$ws.Range("K136").Value = 8931.6312
$ws.Range("M136").Value = -6381.6312
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4510.4375
$ws.Range("J2").Value = 6486.5454
$ws.Range("L2").Value = 38919.2724
$ws.Range("N2").Value = -39145.2724
$ws.Range("H5").Value = 363.83334
$ws.Range("J5").Value = 398.41177
$ws.Range("L5").Value = 1195.23531
$ws.Range("N5").Value = -1419.23531
$ws.Range("H57").Value = 1392.75
$ws.Range("I57").Value = 690.3333
$ws.Range("J57").Value = 3500
$ws.Range("K57").Value = 2070.9999
$ws.Range("L57").Value = 10500
$ws.Range("M57").Value = -1511.9999
$ws.Range("N57").Value = -11618
$ws.Range("H68").Value = 48045.184
$ws.Range("J68").Value = 2896.9443
$ws.Range("L68").Value = 8690.832900000001
$ws.Range("N68").Value = -10312.8329
$ws.Range("H71").Value = 48045.184
$ws.Range("J71").Value = 2896.9443
$ws.Range("L71").Value = 26072.4987
$ws.Range("N71").Value = -34184.4987
$ws.Range("H75").Value = 579.625
$ws.Range("I75").Value = 717.6667
$ws.Range("K75").Value = 2153.0001
$ws.Range("M75").Value = -1155.0001
$ws.Range("H78").Value = 579.625
$ws.Range("I78").Value = 717.6667
$ws.Range("K78").Value = 6459.0003
$ws.Range("M78").Value = -1467.0003
$ws.Range("H97").Value = 380.72726
$ws.Range("I97").Value = 366.33334
$ws.Range("J97").Value = 398
$ws.Range("K97").Value = 1099.00002
$ws.Range("L97").Value = 1194
$ws.Range("M97").Value = -603.0000199999999
$ws.Range("N97").Value = -2186
$ws.Range("H98").Value = 945.5
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 960.6667
$ws.Range("K98").Value = 2700
$ws.Range("L98").Value = 2882.0001
$ws.Range("M98").Value = -1202
$ws.Range("N98").Value = -5878.0001
$ws.Range("H114").Value = 735.5714
$ws.Range("I114").Value = 383.33334
$ws.Range("K114").Value = 1150.00002
$ws.Range("M114").Value = 2103.99998
$ws.Range("H121").Value = 55558924
$ws.Range("J121").Value = 66670690
$ws.Range("L121").Value = 200012070
$ws.Range("N121").Value = -200014690
$ws.Range("H131").Value = 10756508
$ws.Range("I131").Value = 111112520
$ws.Range("J131").Value = 4077.8215
$ws.Range("K131").Value = 333337560
$ws.Range("L131").Value = 12233.4645
$ws.Range("M131").Value = -333332520
$ws.Range("N131").Value = -22313.4645
$ws.Range("H135").Value = 363.83334
$ws.Range("J135").Value = 398.41177
$ws.Range("L135").Value = 3585.70593
$ws.Range("N135").Value = -8655.70593
$ws.Range("H137").Value = 1891.85
$ws.Range("I137").Value = 1077.5
$ws.Range("J137").Value = 2706.2
$ws.Range("K137").Value = 3232.5
$ws.Range("L137").Value = 8118.599999999999
$ws.Range("M137").Value = 1867.5
$ws.Range("N137").Value = -18318.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 38584.6
$ws.Range("J93").Value = 38584.6
$ws.Range("L93").Value = 38584.6
$ws.Range("N93").Value = -42328.6
$ws.Range("H113").Value = 11201
$ws.Range("I113").Value = 1768
$ws.Range("J113").Value = 39500
$ws.Range("K113").Value = 1768
$ws.Range("L113").Value = 39500
$ws.Range("M113").Value = 402
$ws.Range("N113").Value = -43840
$ws.Range("H120").Value = 80000
$ws.Range("J120").Value = 80000
$ws.Range("L120").Value = 80000
$ws.Range("N120").Value = -89676
$ws.Range("H122").Value = 5260.087
$ws.Range("J122").Value = 2955.4443
$ws.Range("L122").Value = 8866.332900000001
$ws.Range("N122").Value = -13766.3329
$ws.Range("H132").Value = 10627.5
$ws.Range("I132").Value = 4953
$ws.Range("J132").Value = 12897.3
$ws.Range("K132").Value = 14859
$ws.Range("L132").Value = 38691.89999999999
$ws.Range("M132").Value = -12329
$ws.Range("N132").Value = -43751.89999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 18523484
$ws.Range("I40").Value = 22732548
$ws.Range("J40").Value = 3601
$ws.Range("K40").Value = 22732548
$ws.Range("L40").Value = 3601
$ws.Range("M40").Value = -22732412
$ws.Range("N40").Value = -3873
$ws.Range("H132").Value = 4682.4287
$ws.Range("J132").Value = 8666
$ws.Range("L132").Value = 25998
$ws.Range("N132").Value = -31058
$ws.Range("H136").Value = 5162.893
$ws.Range("I136").Value = 3086.9443
$ws.Range("K136").Value = 9260.832900000001
$ws.Range("M136").Value = -6710.832900000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3506.5667
$ws.Range("I107").Value = 3508.2693
$ws.Range("K107").Value = 10524.8079
$ws.Range("M107").Value = -8604.8079
